$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109 (pushes existing rows 109:163 down to 110:164,
# inheriting the date number-format from the row above for column D automatically).
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record.
$ws.Cells.Item(109, 1).Value = 10
$ws.Cells.Item(109, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(109, 3).Value = "La Araucanía"
$ws.Cells.Item(109, 4).Value = 44489
$ws.Cells.Item(109, 5).Value = 9
$ws.Cells.Item(109, 6).Value = 100112043
$ws.Cells.Item(109, 7).Value = "Pepino dulce"
$ws.Cells.Item(109, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 75
$ws.Cells.Item(109, 11).Value = 19000
$ws.Cells.Item(109, 12).Value = 20000
$ws.Cells.Item(109, 13).Value = 19467
$ws.Cells.Item(109, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(109, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(109, 16).Value = 1082
$ws.Cells.Item(109, 17).Value = 18
$ws.Cells.Item(109, 18).Value = "Hortaliza"
